# Add new survey wave (28. 9. 2021) columns to both sheets, and bump the
# "aktualizace" date in the title rows, per the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# --- Sheet "data": new column AI, header copies AH1 formatting ---
$ws1.Range("AH1").Copy()
$ws1.Range("AI1").PasteSpecial(-4122)

$ws1.Range("AI1").Value = "28. 9. 2021"
$ws1.Range("AI2").Value = 0.15
$ws1.Range("AI3").Value = 0.14
$ws1.Range("AI4").Value = 0.17
$ws1.Range("AI5").Value = 0.18
$ws1.Range("AI6").Value = 0.36
$ws1.Range("AI7").Value = 0.1
$ws1.Range("AI8").Value = 0.06
$ws1.Range("AI9").Value = 0.17
$ws1.Range("AI10").Value = 0.16
$ws1.Range("AI11").Value = 0.51
$ws1.Range("AI12").Value = 0.08
$ws1.Range("AI13").Value = 0.1
$ws1.Range("AI14").Value = 0.18
$ws1.Range("AI15").Value = 0.18
$ws1.Range("AI16").Value = 0.46
$ws1.Range("AI17").Value = 0.24
$ws1.Range("AI18").Value = 0.22
$ws1.Range("AI19").Value = 0.15
$ws1.Range("AI20").Value = 0.18
$ws1.Range("AI21").Value = 0.21
$ws1.Range("AI22").Value = 0.12
$ws1.Range("AI23").Value = 0.07000000000000001
$ws1.Range("AI24").Value = 0.16
$ws1.Range("AI25").Value = 0.14
$ws1.Range("AI26").Value = 0.51
$ws1.Range("AI27").Value = 0.17
$ws1.Range("AI28").Value = 0.15
$ws1.Range("AI29").Value = 0.18
$ws1.Range("AI30").Value = 0.15
$ws1.Range("AI31").Value = 0.35
$ws1.Range("AI32").Value = 0.16
$ws1.Range("AI33").Value = 0.19
$ws1.Range("AI34").Value = 0.16
$ws1.Range("AI35").Value = 0.26
$ws1.Range("AI36").Value = 0.23
$ws1.Range("AI37").Value = 0.1
$ws1.Range("AI38").Value = 0.11
$ws1.Range("AI39").Value = 0.17
$ws1.Range("AI40").Value = 0.2
$ws1.Range("AI41").Value = 0.42
$ws1.Range("AI42").Value = 0.13
$ws1.Range("AI43").Value = 0.14
$ws1.Range("AI44").Value = 0.18
$ws1.Range("AI45").Value = 0.24
$ws1.Range("AI46").Value = 0.31
$ws1.Range("AI47").Value = 0.17
$ws1.Range("AI48").Value = 0.14
$ws1.Range("AI49").Value = 0.16
$ws1.Range("AI50").Value = 0.15
$ws1.Range("AI51").Value = 0.38
$ws1.Range("AI52").Value = 0.14
$ws1.Range("AI53").Value = 0.11
$ws1.Range("AI54").Value = 0.17
$ws1.Range("AI55").Value = 0.19
$ws1.Range("AI56").Value = 0.39
$ws1.Range("AI57").Value = 0.16
$ws1.Range("AI58").Value = 0.17
$ws1.Range("AI59").Value = 0.17
$ws1.Range("AI60").Value = 0.16
$ws1.Range("AI61").Value = 0.34
$ws1.Range("AI62").Value = 0.17
$ws1.Range("AI63").Value = 0.13
$ws1.Range("AI64").Value = 0.18
$ws1.Range("AI65").Value = 0.17
$ws1.Range("AI66").Value = 0.35
$ws1.Range("AI67").Value = 0.12
$ws1.Range("AI68").Value = 0.14
$ws1.Range("AI69").Value = 0.16
$ws1.Range("AI70").Value = 0.17
$ws1.Range("AI71").Value = 0.41
$ws1.Range("AI72").Value = 0.13
$ws1.Range("AI73").Value = 0.13
$ws1.Range("AI74").Value = 0.16
$ws1.Range("AI75").Value = 0.19
$ws1.Range("AI76").Value = 0.39
$ws1.Range("AI77").Value = 0.14
$ws1.Range("AI78").Value = 0.16
$ws1.Range("AI79").Value = 0.14
$ws1.Range("AI80").Value = 0.21
$ws1.Range("AI81").Value = 0.35
$ws1.Range("AI82").Value = 0.08
$ws1.Range("AI83").Value = 0.04
$ws1.Range("AI84").Value = 0.16
$ws1.Range("AI85").Value = 0.17
$ws1.Range("AI86").Value = 0.55
$ws1.Range("AI87").Value = 0.06
$ws1.Range("AI88").Value = 0.06
$ws1.Range("AI89").Value = 0.2
$ws1.Range("AI90").Value = 0.22
$ws1.Range("AI91").Value = 0.46
$ws1.Range("AI92").Value = 0.24
$ws1.Range("AI93").Value = 0.2
$ws1.Range("AI94").Value = 0.14
$ws1.Range("AI95").Value = 0.17
$ws1.Range("AI96").Value = 0.25
$ws1.Range("AI97").Value = 0.11
$ws1.Range("AI98").Value = 0.09
$ws1.Range("AI99").Value = 0.19
$ws1.Range("AI100").Value = 0.15
$ws1.Range("AI101").Value = 0.46
$ws1.Range("AI102").Value = 0.1
$ws1.Range("AI103").Value = 0.14
$ws1.Range("AI104").Value = 0.15
$ws1.Range("AI105").Value = 0.14
$ws1.Range("AI106").Value = 0.47
$ws1.Range("AI107").Value = 0.23
$ws1.Range("AI108").Value = 0.23
$ws1.Range("AI109").Value = 0.17
$ws1.Range("AI110").Value = 0.18
$ws1.Range("AI111").Value = 0.19

# Title row 112 gets the refreshed "aktualizace" date
$ws1.Range("A112").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# --- Sheet "pocetR": new column AH, header copies AG1 formatting ---
$ws2.Range("AG1").Copy()
$ws2.Range("AH1").PasteSpecial(-4122)

$ws2.Range("AH1").Value = "28. 9. 2021"
$ws2.Range("AH2").Value = 1583
$ws2.Range("AH3").Value = 366
$ws2.Range("AH4").Value = 558
$ws2.Range("AH5").Value = 659
$ws2.Range("AH6").Value = 482
$ws2.Range("AH7").Value = 667
$ws2.Range("AH8").Value = 434
$ws2.Range("AH9").Value = 281
$ws2.Range("AH10").Value = 288
$ws2.Range("AH11").Value = 1014
$ws2.Range("AH12").Value = 794
$ws2.Range("AH13").Value = 789
$ws2.Range("AH14").Value = 821
$ws2.Range("AH15").Value = 367
$ws2.Range("AH16").Value = 191
$ws2.Range("AH17").Value = 204
$ws2.Range("AH18").Value = 195
$ws2.Range("AH19").Value = 302
$ws2.Range("AH20").Value = 298
$ws2.Range("AH21").Value = 171
$ws2.Range("AH22").Value = 256
$ws2.Range("AH23").Value = 361

# Blank trailing placeholder cell AH24 (matches the empty inlineStr cells
# already present across row 24), created via format copy from AG24.
$ws2.Range("AG24").Copy()
$ws2.Range("AH24").PasteSpecial(-4122)

# Title row 24 gets the refreshed "aktualizace" date
$ws2.Range("A24").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"

$excel.CutCopyMode = 0
